$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the gross expenditures and total labor cost inputs
$ws.Range("D3").Value = 253317.02
$ws.Range("D5").Value = 100345.34

# Move the active selection to D3 (cursor moved there after data entry)
$ws.Range("D3").Select()
